# Update cryptocurrency price (column D) and 1h volume change (column E)
# values for the cryptos worksheet, matching the latest scrape.
#
# Each target cell gets its number format pinned to text ("@") before the
# new value is written so that numeric-looking strings (e.g. "1.002",
# "288.35") are preserved verbatim as text instead of being parsed into
# floating point numbers by Excel's type inference. ClearFormats() then
# restores the cell to its original (unstyled/General) appearance, leaving
# only the text value changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "22.450.56"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.568.09"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +0.07%  "
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.17%  "
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "288.35"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -0.57%  "
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3727"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +1.14%  "
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "48.22"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -3.33%  "
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3320"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -2.08%  "
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.133"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -0.78%  "
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07474"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.62%  "
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "20.67"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -2.07%  "
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.933"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -1.09%  "
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.913"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -0.97%  "
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.563.64"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -0.41%  "
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001116"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.27%  "
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.06751"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -0.59%  "
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -2.65%  "
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.347"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "16.44"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +0.82%  "
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.54%  "
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "22.434.37"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +0.11%  "
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.380"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +0.77%  "
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.563"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -3.06%  "
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "152.78"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +2.29%  "
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "19.68"
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.66%  "
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.014"
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.70%  "
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "124.11"
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.54%  "
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.739.99"
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.49%  "
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.055"
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -0.67%  "
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -0.34%  "
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.126"
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -1.34%  "
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "9.754"
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -0.53%  "
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.08300"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -0.86%  "
$cell.ClearFormats()
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.02456"
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -0.77%  "
$cell.ClearFormats()
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2273"
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -0.67%  "
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -1.52%  "
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -0.48%  "
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.286"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -4.39%  "
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +0.53%  "
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.6284"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +1.30%  "
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.84"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -1.46%  "
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.6141"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +5.06%  "
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -0.35%  "
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.048"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.41%  "
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "125.53"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -0.69%  "
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.210"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -1.74%  "
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.07223"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -0.85%  "
$cell.ClearFormats()
